$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.05741433333333334
$ws.Range("H2").Value = 0.172243
$ws.Range("I2").Value = 0.1699029269866134
$ws.Range("J2").Value = 0.1699029269866134
$ws.Range("M2").Value = 0.09920366666666665
$ws.Range("N2").Value = 0.297611
$ws.Range("O2").Value = 0.641640866873065
$ws.Range("P2").Value = 0.6416408668730651
$ws.Range("Q2").Value = 0.005695712385888888
$ws.Range("R2").Value = 0.05126141147299999
$ws.Range("S2").Value = 0.1090166613559617
$ws.Range("T2").Value = 0.1090166613559617

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.05741433333333334
$ws.Range("H3").Value = 0.172243
$ws.Range("I3").Value = 0.1699029269866134
$ws.Range("J3").Value = 0.1699029269866134
$ws.Range("O3").Value = 0.178081099028088
$ws.Range("P3").Value = 0.1780810990280881
$ws.Range("Q3").Value = 0.001580788839666667
$ws.Range("R3").Value = 0.014227099557
$ws.Range("S3").Value = 0.0302564999658651
$ws.Range("T3").Value = 0.03025649996586511

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.05741433333333334
$ws.Range("H4").Value = 0.172243
$ws.Range("I4").Value = 0.1699029269866134
$ws.Range("J4").Value = 0.1699029269866134
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.02787266666666667
$ws.Range("N4").Value = 0.083618
$ws.Range("O4").Value = 0.180278034098847
$ws.Range("P4").Value = 0.180278034098847
$ws.Range("Q4").Value = 0.001600290574888889
$ws.Range("R4").Value = 0.014402615174
$ws.Range("S4").Value = 0.0306297656647866
$ws.Range("T4").Value = 0.0306297656647866

$ws.Range("I5").Value = 0.3444311497741604
$ws.Range("J5").Value = 0.3444311497741605
$ws.Range("M5").Value = 0.09920366666666665
$ws.Range("N5").Value = 0.297611
$ws.Range("O5").Value = 0.641640866873065
$ws.Range("P5").Value = 0.6416408668730651
$ws.Range("Q5").Value = 0.01154648010277778
$ws.Range("R5").Value = 0.103918320925
$ws.Range("S5").Value = 0.2210011015191788
$ws.Range("T5").Value = 0.2210011015191789

$ws.Range("I6").Value = 0.3444311497741604
$ws.Range("J6").Value = 0.3444311497741605
$ws.Range("O6").Value = 0.178081099028088
$ws.Range("P6").Value = 0.1780810990280881
$ws.Range("Q6").Value = 0.003204611758333334
$ws.Range("S6").Value = 0.06133667769129048
$ws.Range("T6").Value = 0.0613366776912905

$ws.Range("I7").Value = 0.3444311497741604
$ws.Range("J7").Value = 0.3444311497741605
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.02787266666666667
$ws.Range("N7").Value = 0.083618
$ws.Range("O7").Value = 0.180278034098847
$ws.Range("P7").Value = 0.180278034098847
$ws.Range("Q7").Value = 0.003244146127777778
$ws.Range("R7").Value = 0.02919731515
$ws.Range("S7").Value = 0.06209337056369118
$ws.Range("T7").Value = 0.06209337056369119

$ws.Range("G8").Value = 0.1641183333333333
$ws.Range("H8").Value = 0.492355
$ws.Range("I8").Value = 0.4856659232392261
$ws.Range("J8").Value = 0.4856659232392261
$ws.Range("M8").Value = 0.09920366666666665
$ws.Range("N8").Value = 0.297611
$ws.Range("O8").Value = 0.641640866873065
$ws.Range("P8").Value = 0.6416408668730651
$ws.Range("Q8").Value = 0.01628114043388889
$ws.Range("R8").Value = 0.146530263905
$ws.Range("S8").Value = 0.3116231039979245
$ws.Range("T8").Value = 0.3116231039979245

$ws.Range("G9").Value = 0.1641183333333333
$ws.Range("H9").Value = 0.492355
$ws.Range("I9").Value = 0.4856659232392261
$ws.Range("J9").Value = 0.4856659232392261
$ws.Range("O9").Value = 0.178081099028088
$ws.Range("P9").Value = 0.1780810990280881
$ws.Range("Q9").Value = 0.004518670071666667
$ws.Range("R9").Value = 0.04066803064500001
$ws.Range("S9").Value = 0.08648792137093242
$ws.Range("T9").Value = 0.08648792137093245

$ws.Range("G10").Value = 0.1641183333333333
$ws.Range("H10").Value = 0.492355
$ws.Range("I10").Value = 0.4856659232392261
$ws.Range("J10").Value = 0.4856659232392261
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.02787266666666667
$ws.Range("N10").Value = 0.083618
$ws.Range("O10").Value = 0.180278034098847
$ws.Range("P10").Value = 0.180278034098847
$ws.Range("Q10").Value = 0.004574415598888889
$ws.Range("R10").Value = 0.04116974039
$ws.Range("S10").Value = 0.08755489787036921
$ws.Range("T10").Value = 0.08755489787036923
